$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Seed the new row's formatting from the row above it (date/border/wrap
# styles already used throughout the table) so no new style entries are
# created, then overwrite with the actual values for the new interview entry.
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A9").Value = 45792
$ws.Range("B9").Value = "Carvale - screening round"
$ws.Range("C9").Value = "5 positive points professionally`n5 negative points professionally`nhow much my manager will rate me out of 10`nwhy you want to switch`nwhere do you want to see yourself in 5 years"

$ws.Rows.Item(9).RowHeight = 72

$ws.Range("C9").Select()
